$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 126 (existing rows 126-167 shift down to 127-168)
$ws.Rows.Item(126).Insert()

# Populate the newly inserted row 126 with the new record
$ws.Cells.Item(126, 1).Value = 9
$ws.Cells.Item(126, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(126, 3).Value = "Metropolitana"
$ws.Cells.Item(126, 4).Value = 44524
$ws.Cells.Item(126, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 5).Value = 13
$ws.Cells.Item(126, 6).Value = 100112030
$ws.Cells.Item(126, 7).Value = "Poroto granado"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 25
$ws.Cells.Item(126, 11).Value = 25000
$ws.Cells.Item(126, 12).Value = 27000
$ws.Cells.Item(126, 13).Value = 25960
$ws.Cells.Item(126, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(126, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(126, 16).Value = 1038
$ws.Cells.Item(126, 17).Value = 25
$ws.Cells.Item(126, 18).Value = "Hortaliza"
